# Scheduled runner update for Brynhildr_Profits price/profit columns (H-N).
# Values below were recalculated from refreshed Market Board averages;
# this script rewrites only the affected leve-profit cells per worksheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 132.42857
$ws.Range("I33").Value = 132.70589
$ws.Range("J33").Value = 131.25
$ws.Range("K33").Value = 132.70589
$ws.Range("L33").Value = 131.25
$ws.Range("M33").Value = 96.29410999999999
$ws.Range("N33").Value = -589.25
# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 4312.6665
$ws.Range("I64").Value = 4002
$ws.Range("J64").Value = 4340.909
$ws.Range("K64").Value = 4002
$ws.Range("L64").Value = 4340.909
$ws.Range("M64").Value = -3754
$ws.Range("N64").Value = -4836.909
# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 4312.6665
$ws.Range("I67").Value = 4002
$ws.Range("J67").Value = 4340.909
$ws.Range("K67").Value = 4002
$ws.Range("L67").Value = 4340.909
$ws.Range("M67").Value = -3144
$ws.Range("N67").Value = -6056.909
# Row 98: The Dotted Line / Enchanted Durium Ink
$ws.Range("H98").Value = 801.55554
$ws.Range("I98").Value = 801.55554
$ws.Range("K98").Value = 801.55554
$ws.Range("M98").Value = 696.44446
# Row 116: Growing Up / Growth Formula Kappa
$ws.Range("H116").Value = 15479.85
$ws.Range("I116").Value = 13849.875
$ws.Range("J116").Value = 16566.5
$ws.Range("K116").Value = 13849.875
$ws.Range("L116").Value = 16566.5
$ws.Range("M116").Value = -10407.875
$ws.Range("N116").Value = -23450.5
# Row 122: Wishful Inking / Enchanted High Durium Ink
$ws.Range("H122").Value = 801.55554
$ws.Range("I122").Value = 801.55554
$ws.Range("K122").Value = 2404.66662
$ws.Range("M122").Value = 45.33338000000003
# Row 125: Body over Mind / Grade 5 Dexterity Alkahest
$ws.Range("H125").Value = 1143.25
$ws.Range("J125").Value = 1241.3334
$ws.Range("L125").Value = 11172.0006
$ws.Range("N125").Value = -16092.0006
# Row 132: Fast-forwarding Flora / Growth Formula Lambda
$ws.Range("H132").Value = 10374.826
$ws.Range("I132").Value = 10814.682
$ws.Range("J132").Value = 698
$ws.Range("K132").Value = 32444.046
$ws.Range("L132").Value = 2094
$ws.Range("M132").Value = -29914.046
$ws.Range("N132").Value = -7154

$ws = $wb.Worksheets.Item("ARM")
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 1596828
$ws.Range("I74").Value = 2320712.8
$ws.Range("J74").Value = 17443.182
$ws.Range("K74").Value = 2320712.8
$ws.Range("L74").Value = 17443.182
$ws.Range("M74").Value = -2319838.8
$ws.Range("N74").Value = -19191.182
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 1596828
$ws.Range("I77").Value = 2320712.8
$ws.Range("J77").Value = 17443.182
$ws.Range("K77").Value = 11603564
$ws.Range("L77").Value = 87215.91
$ws.Range("M77").Value = -11599196
$ws.Range("N77").Value = -95951.91
# Row 132: Don't Bore Me, Ore Me / Mountain Chromite Ingot
$ws.Range("H132").Value = 7716.5
$ws.Range("I132").Value = 5099.8335
$ws.Range("J132").Value = 10333.167
$ws.Range("K132").Value = 15299.5005
$ws.Range("L132").Value = 30999.501
$ws.Range("M132").Value = -12769.5005
$ws.Range("N132").Value = -36059.501

$ws = $wb.Worksheets.Item("BSM")
# Row 15: Anutha Spatha / Bronze Spatha
$ws.Range("H15").Value = 7253.5
$ws.Range("J15").Value = 7253.5
$ws.Range("L15").Value = 7253.5
$ws.Range("N15").Value = -7707.5
# Row 103: The Bigger the Blade / Doman Steel Tachi
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found / Walnut Lumber
$ws.Range("H31").Value = 2394230
$ws.Range("I31").Value = 5850011
$ws.Range("J31").Value = 1766.1538
$ws.Range("K31").Value = 5850011
$ws.Range("L31").Value = 1766.1538
$ws.Range("M31").Value = -5849716
$ws.Range("N31").Value = -2356.1538
# Row 34: Armoires of the Rich and Famous / Walnut Lumber
$ws.Range("H34").Value = 2394230
$ws.Range("I34").Value = 5850011
$ws.Range("J34").Value = 1766.1538
$ws.Range("K34").Value = 5850011
$ws.Range("L34").Value = 1766.1538
$ws.Range("M34").Value = -5849809
$ws.Range("N34").Value = -2170.1538
# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 4316.1665
$ws.Range("I62").Value = 3900
$ws.Range("J62").Value = 4524.25
$ws.Range("K62").Value = 3900
$ws.Range("L62").Value = 4524.25
$ws.Range("M62").Value = -3276
$ws.Range("N62").Value = -5772.25
# Row 64: Almost as Fun as Slingshotting Birds / Cedar Longbow
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 4316.1665
$ws.Range("I65").Value = 3900
$ws.Range("J65").Value = 4524.25
$ws.Range("K65").Value = 19500
$ws.Range("L65").Value = 22621.25
$ws.Range("M65").Value = -16380
$ws.Range("N65").Value = -28861.25
# Row 67: Living Bow to Mouth (L) / Cedar Longbow
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
# Row 68: Do You Even String Bow / Holy Cedar Composite Bow
$ws.Range("H68").Value = 39980
$ws.Range("J68").Value = 39980
$ws.Range("L68").Value = 39980
$ws.Range("N68").Value = -41478
# Row 71: Win One Bow, Get Three Free (L) / Holy Cedar Composite Bow
$ws.Range("H71").Value = 39980
$ws.Range("J71").Value = 39980
$ws.Range("L71").Value = 119940
$ws.Range("N71").Value = -127428
# Row 122: Timber of Tenkonto / Horse Chestnut Lumber
$ws.Range("H122").Value = 6984.6123
$ws.Range("I122").Value = 1887.561
$ws.Range("K122").Value = 5662.683
$ws.Range("M122").Value = -3212.683
# Row 134: Wood You Be Quiet / Ceiba Lumber
$ws.Range("H134").Value = 3173.6155
$ws.Range("I134").Value = 3203.182
$ws.Range("K134").Value = 9609.545999999998
$ws.Range("M134").Value = -7074.545999999998

$ws = $wb.Worksheets.Item("CUL")
# Row 8: Whip It / Sweet Cream
$ws.Range("H8").Value = 244.125
$ws.Range("I8").Value = 244.125
$ws.Range("K8").Value = 732.375
$ws.Range("M8").Value = -593.375
# Row 14: Keep Your Powder Dry / Kukuru Powder
$ws.Range("H14").Value = 319
$ws.Range("I14").Value = 319
$ws.Range("K14").Value = 957
$ws.Range("M14").Value = -784
# Row 50: Moving Up in the World / Rolanberry Cheese
$ws.Range("H50").Value = 461.28946
$ws.Range("I50").Value = 407.05884
$ws.Range("J50").Value = 505.1905
$ws.Range("K50").Value = 1221.17652
$ws.Range("L50").Value = 1515.5715
$ws.Range("M50").Value = -740.17652
$ws.Range("N50").Value = -2477.5715
# Row 53: Rolanberry Fields Forever / Rolanberry Cheese
$ws.Range("H53").Value = 461.28946
$ws.Range("I53").Value = 407.05884
$ws.Range("J53").Value = 505.1905
$ws.Range("K53").Value = 1221.17652
$ws.Range("L53").Value = 1515.5715
$ws.Range("M53").Value = -740.17652
$ws.Range("N53").Value = -2477.5715
# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 2796.8293
$ws.Range("I131").Value = 942.8889
$ws.Range("K131").Value = 2828.6667
$ws.Range("M131").Value = 2211.3333
# Row 133: Friends Are Food / Boiled Alpaca Steak
$ws.Range("H133").Value = 10353
$ws.Range("I133").Value = 7331.636
$ws.Range("J133").Value = 17000
$ws.Range("K133").Value = 21994.908
$ws.Range("L133").Value = 51000
$ws.Range("M133").Value = -16934.908
$ws.Range("N133").Value = -61120

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 29725.354
$ws.Range("I70").Value = 29333.23
$ws.Range("J70").Value = 30999.75
$ws.Range("K70").Value = 29333.23
$ws.Range("L70").Value = 30999.75
$ws.Range("M70").Value = -29063.23
$ws.Range("N70").Value = -31539.75
# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 29725.354
$ws.Range("I73").Value = 29333.23
$ws.Range("J73").Value = 30999.75
$ws.Range("K73").Value = 29333.23
$ws.Range("L73").Value = 30999.75
$ws.Range("M73").Value = -28397.23
$ws.Range("N73").Value = -32871.75
# Row 126: Gold Rush Order / Phrygian Gold Ingot
$ws.Range("H126").Value = 3376.6428
$ws.Range("I126").Value = 3164.75
$ws.Range("J126").Value = 3461.4
$ws.Range("K126").Value = 9494.25
$ws.Range("L126").Value = 10384.2
$ws.Range("M126").Value = -7024.25
$ws.Range("N126").Value = -15324.2
# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 13294.4
$ws.Range("I132").Value = 14436.667
$ws.Range("K132").Value = 43310.001
$ws.Range("M132").Value = -40780.001

$ws = $wb.Worksheets.Item("LTW")
# Row 34: Breeches Served Cold / Goatskin Breeches
$ws.Range("H34").Value = 9000
$ws.Range("I34").Value = 9000
$ws.Range("J34").Value = 9000
$ws.Range("K34").Value = 9000
$ws.Range("L34").Value = 9000
$ws.Range("M34").Value = -8828
$ws.Range("N34").Value = -9344
# Row 68: You Could Say It's a Moving Target / Wyvern Leather
$ws.Range("H68").Value = 14885.883
$ws.Range("I68").Value = 10389.546
$ws.Range("J68").Value = 23129.166
$ws.Range("K68").Value = 10389.546
$ws.Range("L68").Value = 23129.166
$ws.Range("M68").Value = -9640.546
$ws.Range("N68").Value = -24627.166
# Row 71: They Call It Bloody Mary (L) / Wyvern Leather
$ws.Range("H71").Value = 14885.883
$ws.Range("I71").Value = 10389.546
$ws.Range("J71").Value = 23129.166
$ws.Range("K71").Value = 51947.73
$ws.Range("L71").Value = 115645.83
$ws.Range("M71").Value = -48203.73
$ws.Range("N71").Value = -123133.83
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 11908244
$ws.Range("I136").Value = 8931907
$ws.Range("J136").Value = 17860918
$ws.Range("K136").Value = 26795721
$ws.Range("L136").Value = 53582754
$ws.Range("M136").Value = -26793171
$ws.Range("N136").Value = -53587854

$ws = $wb.Worksheets.Item("WVR")
# Row 96: Skills on Display / Ruby Cotton Cloth
$ws.Range("H96").Value = 1309.6154
$ws.Range("I96").Value = 1403
$ws.Range("K96").Value = 1403
$ws.Range("M96").Value = -30
# Row 122: Heavy Armoire / Dark Hempen Cloth
$ws.Range("H122").Value = 25910.604
$ws.Range("I122").Value = 2212.513
$ws.Range("K122").Value = 6637.539
$ws.Range("M122").Value = -4187.539
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 6668697
$ws.Range("I132").Value = 8335026.5
$ws.Range("K132").Value = 25005079.5
$ws.Range("M132").Value = -25002549.5
# Row 140: Glamorous Gloves / Thunderyards Silk Gloves of Casting
$ws.Range("H140").Value = 69984
$ws.Range("J140").Value = 69984
$ws.Range("L140").Value = 69984
$ws.Range("N140").Value = -80344
